$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capacity (column C) and K (financing) values for rows 2-10
$ws.Range("C2:C4").Value = 1000
$ws.Range("K2:K4").Value = 53

$ws.Range("C5:C7").Value = 3000
$ws.Range("K5:K7").Value = 28

$ws.Range("C8:C10").Value = 5000
$ws.Range("K8:K10").Value = 5

# Update the selected range/active cell shown in the saved view
$ws.Range("K2:K4").Select()
